# #79 amend boarding pass background and improve BoardingPassCard layout.
#
# 1) The slide master and every slide layout carry a cached
#    datetimeFigureOut field ("2022/5/19") in their date placeholder —
#    bump the cached text forward a day to "2022/5/20".
# 2) Nudge the "Terminal" / gate-number text boxes (文本框 39 / 文本框 40)
#    that sit at x=6872038 EMU so the boarding-pass card layout lines up
#    with the amended background.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "2022/5/19") {
                $shp.TextFrame.TextRange.Text = "2022/5/20"
            }
        }
    }
}

# --- slide master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- every slide layout off the master ---
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# --- EMU -> point helper (COM geometry is expressed in points) ---
function Emu($v) { return $v / 914400.0 * 72.0 }

# Slide 2: top-level "文本框 39" / "文本框 40" text boxes (boarding pass card)
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.Name -eq "文本框 39" -and [Math]::Round($shp.Left) -eq [Math]::Round((Emu 6872038))) {
        $shp.Top = Emu 1263505
    }
    elseif ($shp.Name -eq "文本框 40" -and [Math]::Round($shp.Left) -eq [Math]::Round((Emu 6872038))) {
        $shp.Top = Emu 1999736
    }
}

# Slide 3: same pair, but nested in the "组合 1" group — only 文本框 40 moves.
$s3 = $p.Slides.Item(3)
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $top = $s3.Shapes.Item($i)
    if ($top.Type -eq 6) {
        for ($k = 1; $k -le $top.GroupItems.Count; $k++) {
            $shp = $top.GroupItems.Item($k)
            if ($shp.Name -eq "文本框 40" -and [Math]::Round($shp.Left) -eq [Math]::Round((Emu 6872038))) {
                $shp.Top = Emu 2002237
            }
        }
    }
}
